# Add year filtering and refresh dataset labels
#
# Adds four new columns (I:L) to Sheet1:
#   I = fed_act_obl_per_1000
#   J = fed_act_obl_indirect_per_1000
#   K = subaward_amount_net_inflow_per_1000
#   L = Year  (populated with "Ave 2020-2024" for every data row)
# and fills in the per-state values for columns I/J/K.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- header row (row 1), matching the bold/bordered/centered style already
# ---- used by the existing headers (copy format from H1, the last header) ----
$ws.Cells.Item(1, 9).Value  = "fed_act_obl_per_1000"
$ws.Cells.Item(1, 10).Value = "fed_act_obl_indirect_per_1000"
$ws.Cells.Item(1, 11).Value = "subaward_amount_net_inflow_per_1000"
$ws.Cells.Item(1, 12).Value = "Year"

$ws.Range("H1").Copy()
$ws.Range("I1:L1").PasteSpecial(-4122)

# ---- per-state data, in the same row order as the existing sheet (rows 2-52) ----
$data = @(
        @(3025326, 3043555, 18228),
        @(8667579, 8401120, -266459),
        @(3179369, 2982672, -196697),
        @(248522, 444627, 196105),
        @(1544643, 1679842, 135198),
        @(2834362, 2819877, -14485),
        @(6247272, 6427607, 180335),
        @(645453, 732677, 87223),
        @(14493118, 14664071, 170952),
        @(1190323, 1242541, 52217),
        @(792407, 807340, 14933),
        @(1836767, 1840403, 3635),
        @(1698388, 1805608, 107220),
        @(1138035, 1178882, 40847),
        @(832012, 874012, 42000),
        @(887588, 1072928, 185340),
        @(1253367, 1369272, 115904),
        @(2441087, 2403590, -37496),
        @(728475, 742811, 14335),
        @(1828955, 1845519, 16564),
        @(8124956, 8064412, -60543),
        @(3152471, 3246163, 93691),
        @(854520, 866437, 11916),
        @(678269, 739584, 61314),
        @(1968796, 1797518, -171278),
        @(2499449, 2128572, -370877),
        @(1037392, 1087843, 50450),
        @(814882, 860364, 45481),
        @(888526, 888709, 183),
        @(1693185, 2569270, 876085),
        @(1313799, 1333757, 19957),
        @(3155488, 3783832, 628343),
        @(1211969, 1263702, 51733),
        @(621463, 672893, 51430),
        @(1338017, 1370915, 32898),
        @(1222590, 940405, -282185),
        @(1474237, 1459953, -14283),
        @(565264, 626282, 61017),
        @(1420967, 1407597, -13370),
        @(639405, 652316, 12910),
        @(2064080, 2092141, 28061),
        @(661097, 760427, 99330),
        @(1532640, 1581062, 48421),
        @(2310576, 2068265, -242310),
        @(1517723, 1841245, 323521),
        @(1260670, 1439670, 178999),
        @(13962630, 13604910, -357720),
        @(1581700, 1584894, 3194),
        @(402271, 669922, 267650),
        @(1697312, 1738660, 41348),
        @(302852, 458409, 155556)
)

for ($idx = 0; $idx -lt $data.Length; $idx++) {
    $row = $idx + 2
    $vals = $data[$idx]
    $ws.Cells.Item($row, 9).Value  = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
    $ws.Cells.Item($row, 11).Value = $vals[2]
    $ws.Cells.Item($row, 12).Value = "Ave 2020-2024"
}
